$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "F7"
$ws.Range("C2").Value = "F3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.859924666666667
$ws.Range("H2").Value = 8.579774
$ws.Range("I2").Value = 0.4493933135217944
$ws.Range("J2").Value = 0.4493933135217944
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 11.28709566666667
$ws.Range("N2").Value = 33.861287
$ws.Range("O2").Value = 0.2158446779538137
$ws.Range("P2").Value = 0.2158446779538137
$ws.Range("Q2").Value = 32.28024331212644
$ws.Range("R2").Value = 290.5221898091381
$ws.Range("S2").Value = 0.09699915503170892
$ws.Range("T2").Value = 0.09699915503170892

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "F7"
$ws.Range("C3").Value = "F3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.859924666666667
$ws.Range("H3").Value = 8.579774
$ws.Range("I3").Value = 0.4493933135217944
$ws.Range("J3").Value = 0.4493933135217944
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 39.44626233333333
$ws.Range("N3").Value = 118.338787
$ws.Range("O3").Value = 0.75433628289025
$ws.Range("P3").Value = 0.75433628289025
$ws.Range("Q3").Value = 112.8133386549042
$ws.Range("R3").Value = 1015.320047894138
$ws.Range("S3").Value = 0.3389936816777631
$ws.Range("T3").Value = 0.3389936816777631

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "F7"
$ws.Range("C4").Value = "F3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.859924666666667
$ws.Range("H4").Value = 8.579774
$ws.Range("I4").Value = 0.4493933135217944
$ws.Range("J4").Value = 0.4493933135217944
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.126152
$ws.Range("N4").Value = 0.378456
$ws.Range("O4").Value = 0.002412421992102324
$ws.Range("P4").Value = 0.002412421992102323
$ws.Range("Q4").Value = 0.3607852165493334
$ws.Range("R4").Value = 3.247066948944
$ws.Range("S4").Value = 0.001084126312643711
$ws.Range("T4").Value = 0.001084126312643711

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "F7"
$ws.Range("C5").Value = "F3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.859924666666667
$ws.Range("H5").Value = 8.579774
$ws.Range("I5").Value = 0.4493933135217944
$ws.Range("J5").Value = 0.4493933135217944
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.433165333333333
$ws.Range("N5").Value = 4.299496
$ws.Range("O5").Value = 0.02740661716383402
$ws.Range("P5").Value = 0.02740661716383402
$ws.Range("Q5").Value = 4.098744888211555
$ws.Range("R5").Value = 36.888703993904
$ws.Range("S5").Value = 0.01231635049967865
$ws.Range("T5").Value = 0.01231635049967865

$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "F7"
$ws.Range("C6").Value = "F3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.504043333333333
$ws.Range("H6").Value = 10.51213
$ws.Range("I6").Value = 0.5506066864782057
$ws.Range("J6").Value = 0.5506066864782055
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 11.28709566666667
$ws.Range("N6").Value = 33.861287
$ws.Range("O6").Value = 0.2158446779538137
$ws.Range("P6").Value = 0.2158446779538137
$ws.Range("Q6").Value = 39.55047232347889
$ws.Range("R6").Value = 355.95425091131
$ws.Range("S6").Value = 0.1188455229221047
$ws.Range("T6").Value = 0.1188455229221047

$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "F7"
$ws.Range("C7").Value = "F3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.504043333333333
$ws.Range("H7").Value = 10.51213
$ws.Range("I7").Value = 0.5506066864782057
$ws.Range("J7").Value = 0.5506066864782055
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 39.44626233333333
$ws.Range("N7").Value = 118.338787
$ws.Range("O7").Value = 0.75433628289025
$ws.Range("P7").Value = 0.75433628289025
$ws.Range("Q7").Value = 138.2214125540344
$ws.Range("R7").Value = 1243.99271298631
$ws.Range("S7").Value = 0.415342601212487
$ws.Range("T7").Value = 0.4153426012124868

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "F7"
$ws.Range("C8").Value = "F3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.504043333333333
$ws.Range("H8").Value = 10.51213
$ws.Range("I8").Value = 0.5506066864782057
$ws.Range("J8").Value = 0.5506066864782055
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.126152
$ws.Range("N8").Value = 0.378456
$ws.Range("O8").Value = 0.002412421992102324
$ws.Range("P8").Value = 0.002412421992102323
$ws.Range("Q8").Value = 0.4420420745866667
$ws.Range("R8").Value = 3.97837867128
$ws.Range("S8").Value = 0.001328295679458612
$ws.Range("T8").Value = 0.001328295679458612

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "F7"
$ws.Range("C9").Value = "F3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.504043333333333
$ws.Range("H9").Value = 10.51213
$ws.Range("I9").Value = 0.5506066864782057
$ws.Range("J9").Value = 0.5506066864782055
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.433165333333333
$ws.Range("N9").Value = 4.299496
$ws.Range("O9").Value = 0.02740661716383402
$ws.Range("P9").Value = 0.02740661716383402
$ws.Range("Q9").Value = 5.021873431831111
$ws.Range("R9").Value = 45.19686088647999
$ws.Range("S9").Value = 0.01509026666415537
$ws.Range("T9").Value = 0.01509026666415536
